$d = $word.ActiveDocument

function Is-EmptyPara($para) {
    return ($para.Range.Text.TrimEnd([char]13, [char]7).Length -eq 0)
}

# 1. Remove the existing "_GoBack" bookmark from the "Project: ...SHORT~" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Of the two empty "Inside Address" paragraphs just before "Fort Shafter, MP T-118
#    Renovation", drop the second one and move the "_GoBack" bookmark onto the first.
#    (Handled before the NOAA insertion below so paragraph indices used here stay valid.)
$targetPara = $null
for ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ((Is-EmptyPara $para) -and ($para.Style.NameLocal -eq "Inside Address")) {
        $nextPara = $d.Paragraphs.Item($i + 1)
        if ((Is-EmptyPara $nextPara) -and ($nextPara.Style.NameLocal -eq "Inside Address")) {
            $targetPara = $para
            $nextPara.Range.Delete()
            break
        }
    }
}
if ($targetPara -ne $null) {
    $d.Bookmarks.Add("_GoBack", $targetPara.Range) | Out-Null
}

# 3. Insert a new empty paragraph (indented 720 twips / 0.5in) right after the
#    "NOAA Pacific Regional Center, Main Facility" paragraph.
$noaaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "NOAA Pacific Regional Center, Main Facility") {
        $noaaPara = $para
        break
    }
}
$insertPoint = $d.Range($noaaPara.Range.End, $noaaPara.Range.End)
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $insertPoint.InsertXML($xml)
